# Add a new "2023" column (T) to the Hepatitis B incidence table, mirroring
# the formatting of the existing "2022" column (S) and extending row 1's
# height slightly (row 4 gets an explicit custom height too).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column T, keyed by row number (row 3 is the "2023" header;
# rows 4-33 are the data rows).
$values = @{
    3  = 2023
    4  = 2.3381104968484805
    5  = 2.0344672190198714
    6  = 2.6483752218014245
    7  = 3.9852372948902328
    8  = 4.5532396299967433
    9  = 3.4291318466903733
    10 = 1.2089851778417198
    11 = 1.521116134174612
    12 = 0.9008846687447073
    13 = 3.694303753043183
    14 = 4.0607488020791038
    15 = 3.327319511401615
    16 = 0.32236434908190637
    17 = 0
    18 = 0.63756806039044667
    19 = 2.1691385808410835
    20 = 1.5024572004578396
    21 = 2.8259763748375066
    22 = 6.1744985943935555
    23 = 4.3993752887090034
    24 = 7.9169155696940479
    25 = 2.8763040791558883
    26 = 1.4751329463567904
    27 = 4.2954684675262591
    28 = 1.8177568880002077
    29 = 1.581380197008345
    30 = 2.103608453446189
    31 = 1.3736037318066185
    32 = 2.249820014398848
    33 = 0.53701655085009725
}

# Column T should look exactly like column S for each row: same number
# format / style. Copy S's formatting into T first, then write the value.
for ($row = 3; $row -le 33; $row++) {
    $src = $ws.Cells.Item($row, 19)   # column S
    $dst = $ws.Cells.Item($row, 20)   # column T

    $src.Copy()
    $dst.PasteSpecial(-4122)          # xlPasteFormats

    $dst.Value2 = $values[$row]
}

$excel.CutCopyMode = $false

# Row 4 picks up an explicit custom height in the edited workbook.
$ws.Rows.Item(4).RowHeight = 16.5

# Reset the selection away from the old "T3" placeholder cell left over from
# editing, back to the sheet's natural default (A1).
$ws.Range("A1").Select()
